# Auto-generated edit script applying numeric corrections to Phantom_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 240
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("H5").Value = 284.7143
$ws.Range("I5").Value = 328.6
$ws.Range("J5").Value = 175
$ws.Range("K5").Value = 328.6
$ws.Range("L5").Value = 175
$ws.Range("M5").Value = -213.6
$ws.Range("N5").Value = -405
$ws.Range("H32").Value = 3757.6
$ws.Range("I32").Value = 2266.6667
$ws.Range("J32").Value = 5994
$ws.Range("K32").Value = 2266.6667
$ws.Range("L32").Value = 5994
$ws.Range("M32").Value = -1940.6667
$ws.Range("N32").Value = -6646
$ws.Range("H37").Value = 383.33334
$ws.Range("I37").Value = 260
$ws.Range("K37").Value = 780
$ws.Range("M37").Value = -654
$ws.Range("H39").Value = 191
$ws.Range("I39").Value = 95
$ws.Range("J39").Value = 287
$ws.Range("K39").Value = 285
$ws.Range("L39").Value = 861
$ws.Range("M39").Value = 11
$ws.Range("N39").Value = -1453
$ws.Range("H40").Value = 1203.5
$ws.Range("I40").Value = 1221.5385
$ws.Range("K40").Value = 1221.5385
$ws.Range("M40").Value = -1046.5385
$ws.Range("H49").Value = 625
$ws.Range("I49").Value = 625
$ws.Range("K49").Value = 1875
$ws.Range("M49").Value = -1739
$ws.Range("H54").Value = 12500
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 15000
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 15000
$ws.Range("M54").Value = -9514
$ws.Range("N54").Value = -15972
$ws.Range("H64").Value = 4949.6665
$ws.Range("I64").Value = 4939.8
$ws.Range("K64").Value = 4939.8
$ws.Range("M64").Value = -4691.8
$ws.Range("H67").Value = 4949.6665
$ws.Range("I67").Value = 4939.8
$ws.Range("K67").Value = 4939.8
$ws.Range("M67").Value = -4081.8
$ws.Range("H136").Value = 240000
$ws.Range("J136").Value = 240000
$ws.Range("L136").Value = 240000
$ws.Range("N136").Value = -250200
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("H5").Value = 350.83334
$ws.Range("I5").Value = 350.83334
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 350.83334
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -238.83334
$ws.Range("H61").Value = 3133.3809
$ws.Range("I61").Value = 2961.1667
$ws.Range("J61").Value = 4166.6665
$ws.Range("K61").Value = 2961.1667
$ws.Range("L61").Value = 4166.6665
$ws.Range("M61").Value = -2749.1667
$ws.Range("N61").Value = -4590.6665
$ws.Range("H132").Value = 2726.7
$ws.Range("I132").Value = 2640.2354
$ws.Range("J132").Value = 3216.6667
$ws.Range("K132").Value = 7920.706200000001
$ws.Range("L132").Value = 9650.000100000001
$ws.Range("M132").Value = -5390.706200000001
$ws.Range("N132").Value = -14710.0001
$ws.Range("H133").Value = 79947.5
$ws.Range("J133").Value = 79947.5
$ws.Range("L133").Value = 79947.5
$ws.Range("N133").Value = -85007.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 3133.3809
$ws.Range("I136").Value = 2961.1667
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 8883.500100000001
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -6333.500100000001
$ws.Range("N136").Value = -17599.9995
$ws.Range("M4").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 350.83334
$ws.Range("I4").Value = 350.83334
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 350.83334
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -235.83334
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 204.86957
$ws.Range("I7").Value = 157.05263
$ws.Range("K7").Value = 157.05263
$ws.Range("M7").Value = -44.05262999999999
$ws.Range("H16").Value = 839.4
$ws.Range("I16").Value = 839.4
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 839.4
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -552.4
$ws.Range("H31").Value = 3529.8
$ws.Range("J31").Value = 5742
$ws.Range("L31").Value = 5742
$ws.Range("N31").Value = -6332
$ws.Range("H34").Value = 3529.8
$ws.Range("J34").Value = 5742
$ws.Range("L34").Value = 5742
$ws.Range("N34").Value = -6146
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H94").Value = 1461.5
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("H99").Value = 3011.5
$ws.Range("I99").Value = 3011.5
$ws.Range("K99").Value = 3011.5
$ws.Range("M99").Value = -1513.5
$ws.Range("H113").Value = 839.4
$ws.Range("I113").Value = 839.4
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 839.4
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1330.6
$ws.Range("H126").Value = 3011.5
$ws.Range("I126").Value = 3011.5
$ws.Range("K126").Value = 9034.5
$ws.Range("M126").Value = -6564.5
$ws.Range("H135").Value = 79999
$ws.Range("J135").Value = 79997
$ws.Range("L135").Value = 79997
$ws.Range("N135").Value = -90137
$ws.Range("N16").ClearContents()
$ws.Range("N52").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2640.2222
$ws.Range("J113").Value = 2510.5
$ws.Range("L113").Value = 7531.5
$ws.Range("N113").Value = -11871.5
$ws.Range("H140").Value = 716257.9
$ws.Range("I140").Value = 716257.9
$ws.Range("K140").Value = 2148773.7
$ws.Range("M140").Value = -2143593.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 1507.5
$ws.Range("I35").Value = 15
$ws.Range("K35").Value = 15
$ws.Range("M35").Value = 283
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H80").Value = 2593.25
$ws.Range("I80").Value = 2452
$ws.Range("J80").Value = 2734.5
$ws.Range("K80").Value = 2452
$ws.Range("L80").Value = 2734.5
$ws.Range("M80").Value = -1454
$ws.Range("N80").Value = -4730.5
$ws.Range("H83").Value = 2593.25
$ws.Range("I83").Value = 2452
$ws.Range("J83").Value = 2734.5
$ws.Range("K83").Value = 12260
$ws.Range("L83").Value = 13672.5
$ws.Range("M83").Value = -7268
$ws.Range("N83").Value = -23656.5
$ws.Range("N49").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 674.75
$ws.Range("I22").Value = 448.375
$ws.Range("J22").Value = 1127.5
$ws.Range("K22").Value = 448.375
$ws.Range("L22").Value = 1127.5
$ws.Range("M22").Value = -153.375
$ws.Range("N22").Value = -1717.5
$ws.Range("H27").Value = 674.75
$ws.Range("I27").Value = 448.375
$ws.Range("J27").Value = 1127.5
$ws.Range("K27").Value = 448.375
$ws.Range("L27").Value = 1127.5
$ws.Range("M27").Value = -341.375
$ws.Range("N27").Value = -1341.5
$ws.Range("H50").Value = 16000
$ws.Range("I50").Value = 16000
$ws.Range("K50").Value = 16000
$ws.Range("M50").Value = -15363
$ws.Range("H55").Value = 1042.9
$ws.Range("I55").Value = 369.75
$ws.Range("J55").Value = 1491.6666
$ws.Range("K55").Value = 369.75
$ws.Range("L55").Value = 1491.6666
$ws.Range("M55").Value = -196.75
$ws.Range("N55").Value = -1837.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 57497.5
$ws.Range("I2").Value = 44995
$ws.Range("J2").Value = 70000
$ws.Range("K2").Value = 44995
$ws.Range("L2").Value = 70000
$ws.Range("M2").Value = -44883
$ws.Range("N2").Value = -70224
$ws.Range("H69").Value = 75135.5
$ws.Range("J69").Value = 75135.5
$ws.Range("L69").Value = 75135.5
$ws.Range("N69").Value = -76633.5
$ws.Range("H70").Value = 38949.5
$ws.Range("I70").Value = 38949.5
$ws.Range("K70").Value = 38949.5
$ws.Range("M70").Value = -38634.5
$ws.Range("H72").Value = 75135.5
$ws.Range("J72").Value = 75135.5
$ws.Range("L72").Value = 225406.5
$ws.Range("N72").Value = -232894.5
$ws.Range("H73").Value = 38949.5
$ws.Range("I73").Value = 38949.5
$ws.Range("K73").Value = 38949.5
$ws.Range("M73").Value = -37857.5
$ws.Range("H122").Value = 4909
$ws.Range("I122").Value = 5727.143
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 17181.429
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -14731.429
$ws.Range("N122").Value = -13900

Write-Host "Applied 234 cell updates across 8 sheets"
